# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: C2 was an erroneous stray value from the naive forecaster bug - remove it.
$ws.Range("C2").ClearContents()

# Recalculated forecast values (tiny precision corrections from the bug fix).
$ws.Range("E4").Value = 2.828066716168043

$ws.Range("C5").Value = 1.891592186533786
$ws.Range("E5").Value = 2.544631191216373

$ws.Range("E6").Value = 1.312870290004309

$ws.Range("C7").Value = 0.618025493879526
$ws.Range("E7").Value = 0.7749619016294229

$ws.Range("E8").Value = 0.7487574275251818

$ws.Range("C9").Value = 1.905862317202112

$ws.Range("C11").Value = 1.917627847674042
$ws.Range("E11").Value = 2.69471174461664

$ws.Range("C12").Value = 1.119562422009124

$ws.Range("C13").Value = 1.344920716048215
$ws.Range("E13").Value = 1.037735724446587

$ws.Range("C14").Value = 2.195375580740744

$ws.Range("C15").Value = 2.491319804758541
$ws.Range("E15").Value = 2.78106797904647

$ws.Range("E16").Value = 1.872042068954638

$ws.Range("C17").Value = -4.149799191324066
$ws.Range("E17").Value = -2.475922651815632

$ws.Range("C19").Value = 2.221594549876427
